$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Gender" header in column L, row 1
$ws.Range("L1").Value = "Gender"

# Update the selection to match the diff (activeCell M1, sqref M1)
$ws.Range("M1").Select()
